$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.630.82"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.591.65"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.71%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.87"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.512"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.86%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  -1.39%  "

$ws.Range("E9").Value = "  -2.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.58"
$ws.Range("D10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0833"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.815.92"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.600.63"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.75%  "

$ws.Range("E15").Value = "  -2.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.01"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.596.74"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.11"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.65%  "

$ws.Range("E21").Value = "  -2.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.25"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.30"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.87"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.88"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.48%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.115"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.30"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0505"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.62%  "

$ws.Range("E31").Value = "  -1.31%  "

$ws.Range("E32").Value = "  -3.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.658"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -10.28%  "

$ws.Range("E34").Value = "  -3.36%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.293.40"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.65%  "

$ws.Range("E36").Value = "  -0.44%  "

$ws.Range("E37").Value = "  -4.95%  "

$ws.Range("E38").Value = "  -3.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.828"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.18%  "

$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("E42").Value = "  +0.66%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.18"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.97%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.09"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.19%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.728.06"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.02"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.73%  "

$ws.Range("E47").Value = "  -1.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.801"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0982"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.51"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.75%  "
